$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Março")
$ws.Activate()

# New expense entries for March (rows 20-22), continuing the existing list.
$ws.Range("A20").Value = "Tablet Expad 8"
$ws.Range("B20").Value = 670.05
$ws.Range("A21").Value = "Teclado tablet"
$ws.Range("B21").Value = 209.09
$ws.Range("A22").Value = "Pedidso mãe"
$ws.Range("B22").Value = 227

# Match the cursor/selection state left behind after the edit.
$ws.Range("B23").Select()
